# The deck's footer "date last updated" field reads 11/6/16 everywhere
# it appears (the slide master and all 11 slide layouts - the single
# content slide itself doesn't instantiate a footer/date placeholder).
# Bump it to 11/12/16 to match the refreshed export date.

$p = $ppt.ActivePresentation
$oldDate = "11/6/16"
$newDate = "11/12/16"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh) {
            $cur = $shp.TextFrame.TextRange.Text
            if ($cur -eq $oldDate -or $cur -ne $newDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Any slide that overrides the footer date placeholder directly.
$slides = $p.Slides
for ($si = 1; $si -le $slides.Count; $si++) {
    Update-DatePlaceholder $slides.Item($si).Shapes
}
